$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Updated Cypher query for the "FilesTab" row (row 4, column B = "query").
# This replaces the old query (no OPTIONAL MATCH / sample / Format columns)
# with the new one that adds sample linkage, a Format column, a formatted
# Size column and a Sample ID column, and orders/limits the results.
$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN coalesce(f.file_name, '') AS `File Name`, 
 coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`, 
      CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
   coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis 
        Order By f.file_name LIMIT 100
'@

$ws.Range("B4").Value = $newFilesTabQuery

# Row heights reflow slightly because of the new (differently sized) query
# text in row 4 -- row 4 now holds the longest text and grows to the
# maximum Excel row height, while rows 2-3 shrink a touch.
$ws.Rows.Item(2).RowHeight = 348
$ws.Rows.Item(3).RowHeight = 275.5
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# Selection / scroll position moved up one row (from C5 to C4).
$ws.Activate()
$ws.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Window size/position of the workbook's own window, as recorded the
# last time it was saved.
$win = $excel.ActiveWindow
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 10420
